# Weekly fruit/vegetable price update: a new price record is inserted as
# the most recent entry (row 5), pushing all existing records down by one
# row (old row 5 -> row 6, old row 6 -> row 7, ..., old row 36 -> row 37).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5; Excel shifts rows 5:36 down to
# 6:37 and carries the formatting (incl. the date number format on column D)
# down with them.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = 'Vega Modelo de Temuco'
$ws.Range("C5").Value = 'La Araucanía'
$ws.Range("D5").Value = 44532
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100114002
$ws.Range("G5").Value = 'Camote'
$ws.Range("H5").Value = 'Sin especificar'
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 18000
$ws.Range("N5").Value = '$/malla 20 kilos'
$ws.Range("O5").Value = 'Perú'
$ws.Range("P5").Value = 900
$ws.Range("Q5").Value = 20
$ws.Range("R5").Value = 'Hortaliza'
